$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date string (slashes -> dashes)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    # Prefix with an apostrophe so Excel keeps the dashed date as literal
    # text instead of re-interpreting it as a date value/serial number.
    $ws.Range("A$row").Value = "'" + $dates[$row]
}

# Map of row -> D,E,F,G,H values that changed
$values = @{
    3  = @(1,0,0,1,1)
    4  = @(1,1,0,0,0)
    5  = @(1,1,0,0,0)
    6  = @(1,1,0,0,0)
    14 = @(1,1,0,0,0)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("D$row").Value = $v[0]
    $ws.Range("E$row").Value = $v[1]
    $ws.Range("F$row").Value = $v[2]
    $ws.Range("G$row").Value = $v[3]
    $ws.Range("H$row").Value = $v[4]
}
